$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ID column (A2:A4) got re-labelled: C1/C2/C3 -> C4/C3/C1
$ws.Range("A2").Value = "C4"
$ws.Range("A3").Value = "C3"
$ws.Range("A4").Value = "C1"

# E2 / E3 in the R1 column flip from "T" to a literal "-" (entered with a
# leading apostrophe so Excel stores it as text rather than a minus sign)
$ws.Range("E2").Value = "'-"
$ws.Range("E3").Value = "'-"

# Update the active selection to match the author's final cursor position
$ws.Range("E4").Select()
